$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic re-ordering of the species-observation rows
# 25-31: the block of columns A,D,E,F,G,H,Q,R (taxon id / redlist status /
# taxon id / name / scientific name / author / coordinates) moves from one
# row to another, while column B (Taxonsorteringsordning) is independently
# re-assigned new values. In addition, the blank Age/Gender/Activity/Method
# markers (K-N) and the public comment (AC) move from row 29 to row 31.

$cols = @("A","D","E","F","G","H","Q","R")

# Snapshot the "before" values of the columns that move as a block.
$before = @{}
for ($r = 25; $r -le 31; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value()
    }
    $before[$r] = $row
}

# Row r gets the block of data (A,D,E,F,G,H,Q,R) that used to live on row
# $srcMap[r].
$srcMap = @{
    25 = 26
    26 = 31
    27 = 30
    28 = 27
    29 = 28
    30 = 25
    31 = 29
}

foreach ($r in 25..31) {
    $src = $before[$srcMap[$r]]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $src[$c]
    }
}

# Column B (Taxonsorteringsordning) gets new values independent of the
# permutation above.
$newB = @{
    25 = 89557
    26 = 89557
    27 = 77636
    28 = 90099
    29 = 90221
    30 = 56446
    31 = 56430
}
foreach ($r in 25..31) {
    $ws.Range("B$r").Value = $newB[$r]
}

# The blank Age/Gender/Activity/Method markers move from row 29 to row 31.
# (ClearContents on row 29 removes the cells entirely, matching the diff;
# writing a real empty value directly would drop the cell too, so force a
# blank-text cell via the text quote-prefix then strip the resulting
# formatting back to Normal.)
foreach ($c in @("K", "L", "M", "N")) {
    $ws.Range("$c" + "29").ClearContents()
    $cell = $ws.Range("$c" + "31")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# The public comment moves from row 29 to row 31.
$ws.Range("AC29").ClearContents()
$ws.Range("AC31").Value = "ringhack äldre"
